# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) indicating Control (0) vs MDD (1) patients,
# and updates the D8/E8 prediction/error values to reflect the refit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Label" header in H1, matching the formatting of the
#     other header cells (bold, centered, bordered) by copying B1's format.
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 8).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 8).Value = "Label"

# --- Populate the H column: 0 for Control patients, 1 for MDD patients.
$labels = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}

# --- Update the refit prediction/error values for row 8 (MDD 3, 100 iters).
$ws.Range("D8").Value = 0.9137627374429688
$ws.Range("E8").Value = 0.08623726255703124
